$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing currency-rate table occupies A2:H28. Duplicate that whole
# block immediately below itself (starting at row 29), so the sheet ends
# up holding the same data twice (rows 2-28 and rows 29-55). Using
# Range.Copy (rather than assigning .Value/.Value2 arrays) preserves the
# original shared-string cell types instead of Excel's usual
# "numeric-looking text becomes a number" auto-conversion.
$src = $ws.Range("A2:H28")
$dst = $ws.Range("A29")
$src.Copy($dst)
